$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated Price (column D) and Volume/1h (column E) values from the
# latest symbol-list refresh. Cells are stored as literal text (not numbers)
# in the source sheet, so NumberFormat is forced to "@" (Text) before the
# assignment to stop Excel from auto-converting numeric- or percent-looking
# strings (e.g. "306.06", "-0.85%") into real numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "306.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.85%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.05"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "7.29%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.112"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.19%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08071"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.55%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.933"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.56%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.198"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.23%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "8.057"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.49%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9260"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.04%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1394"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-6.39%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1919"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.14%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09026"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.93%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03519"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.32%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09791"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.74%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001405"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.33%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005880"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-9.52%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.765"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.07%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3461"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1325"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.98%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.677"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2415"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "3.08%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04376"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.06%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001206"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.27%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.79%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02041"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-4.01%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05029"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.80%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007528"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.78%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009710"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.57%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.83%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.83%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009799"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.97%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006210"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.01%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.06%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002783"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "12.61%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.06%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.06%"
